# Quarterly update: roll the data window forward by one quarter.
# Drop the oldest quarter column of data, shift all remaining quarters left,
# and append the newly-reported quarter (12 ماهه منتهی به 1401/12, published 1402-02-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: quarter-period header labels (rolled forward by one quarter) ---
$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish-date header labels (rolled forward by one quarter) ---
$ws.Range("D9").Value = "1400-10-29 (3)"
$ws.Range("E9").Value = "1401-03-11 (8)"
$ws.Range("F9").Value = "1401-04-29 (2)"
$ws.Range("G9").Value = "1401-08-29 (4)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-02-30 (7)"
$ws.Range("J9").Value = "1401-04-29"
$ws.Range("K9").Value = "1401-08-29 (2)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-02-30"

# --- Data rows 11-27: values rolled forward by one quarter (oldest dropped, newest appended) ---
# Row 11
$ws.Range("D11").Value = 23137
$ws.Range("E11").Value = 37454
$ws.Range("F11").Value = 10811
$ws.Range("G11").Value = 23702
$ws.Range("H11").Value = 36868
$ws.Range("I11").Value = 51890
$ws.Range("J11").Value = 14707
$ws.Range("K11").Value = 34873
$ws.Range("L11").Value = 54004
$ws.Range("M11").Value = 69840

# Row 12
$ws.Range("D12").Value = -21624
$ws.Range("E12").Value = -32346
$ws.Range("F12").Value = -9171
$ws.Range("G12").Value = -19779
$ws.Range("H12").Value = -31244
$ws.Range("I12").Value = -42900
$ws.Range("J12").Value = -13082
$ws.Range("K12").Value = -26853
$ws.Range("L12").Value = -43389
$ws.Range("M12").Value = -56066

# Row 13
$ws.Range("D13").Value = 1513
$ws.Range("E13").Value = 5108
$ws.Range("F13").Value = 1640
$ws.Range("G13").Value = 3923
$ws.Range("H13").Value = 5624
$ws.Range("I13").Value = 8990
$ws.Range("J13").Value = 1625
$ws.Range("K13").Value = 8020
$ws.Range("L13").Value = 10614
$ws.Range("M13").Value = 13774

# Row 14
$ws.Range("D14").Value = -1907
$ws.Range("E14").Value = -3503
$ws.Range("F14").Value = -841
$ws.Range("G14").Value = -1902
$ws.Range("H14").Value = -2814
$ws.Range("I14").Value = -4144
$ws.Range("J14").Value = -875
$ws.Range("K14").Value = -2238
$ws.Range("L14").Value = -3091
$ws.Range("M14").Value = -3853

# Row 16
$ws.Range("D16").Value = 2083
$ws.Range("E16").Value = 1820
$ws.Range("F16").Value = -584
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 134
$ws.Range("I16").Value = 17
$ws.Range("J16").Value = 77
$ws.Range("K16").Value = 119
$ws.Range("L16").Value = 956
$ws.Range("M16").Value = 2503

# Row 17
$ws.Range("D17").Value = 1688
$ws.Range("E17").Value = 3425
$ws.Range("F17").Value = 215
$ws.Range("G17").Value = 2032
$ws.Range("H17").Value = 2945
$ws.Range("I17").Value = 4864
$ws.Range("J17").Value = 827
$ws.Range("K17").Value = 5901
$ws.Range("L17").Value = 8479
$ws.Range("M17").Value = 12423

# Row 18
$ws.Range("D18").Value = -76
$ws.Range("E18").Value = -76
$ws.Range("F18").Value = -45
$ws.Range("G18").Value = -93
$ws.Range("H18").Value = -126
$ws.Range("I18").Value = -183
$ws.Range("J18").Value = -30
$ws.Range("K18").Value = -45
$ws.Range("L18").Value = -49
$ws.Range("M18").Value = -44

# Row 19
$ws.Range("D19").Value = 496
$ws.Range("E19").Value = 659
$ws.Range("F19").Value = 312
$ws.Range("G19").Value = 330
$ws.Range("H19").Value = 511
$ws.Range("I19").Value = 817
$ws.Range("J19").Value = 92
$ws.Range("K19").Value = 179
$ws.Range("L19").Value = 400
$ws.Range("M19").Value = 400

# Row 20
$ws.Range("D20").Value = 2108
$ws.Range("E20").Value = 4008
$ws.Range("F20").Value = 481
$ws.Range("G20").Value = 2270
$ws.Range("H20").Value = 3330
$ws.Range("I20").Value = 5498
$ws.Range("J20").Value = 890
$ws.Range("K20").Value = 6034
$ws.Range("L20").Value = 8830
$ws.Range("M20").Value = 12779

# Row 21
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = -22
$ws.Range("F21").Value = "-"
$ws.Range("G21").Value = -131
$ws.Range("H21").Value = -256
$ws.Range("I21").Value = -352
$ws.Range("J21").Value = -113
$ws.Range("K21").Value = -871
$ws.Range("L21").Value = -1172
$ws.Range("M21").Value = -1448

# Row 22
$ws.Range("D22").Value = 2108
$ws.Range("E22").Value = 3985
$ws.Range("F22").Value = 481
$ws.Range("G22").Value = 2138
$ws.Range("H22").Value = 3074
$ws.Range("I22").Value = 5146
$ws.Range("J22").Value = 777
$ws.Range("K22").Value = 5164
$ws.Range("L22").Value = 7658
$ws.Range("M22").Value = 11331

# Row 24
$ws.Range("D24").Value = 2108
$ws.Range("E24").Value = 3985
$ws.Range("F24").Value = 481
$ws.Range("G24").Value = 2138
$ws.Range("H24").Value = 3074
$ws.Range("I24").Value = 5146
$ws.Range("J24").Value = 777
$ws.Range("K24").Value = 5164
$ws.Range("L24").Value = 7658
$ws.Range("M24").Value = 11331

# Row 26
$ws.Range("D26").Value = 2560
$ws.Range("E26").Value = 8076
$ws.Range("F26").Value = 7827
$ws.Range("G26").Value = 7369
$ws.Range("H26").Value = 7019
$ws.Range("I26").Value = 6921
$ws.Range("J26").Value = 6214
$ws.Range("K26").Value = 6047
$ws.Range("L26").Value = 5742
$ws.Range("M26").Value = 5174

